$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("English and Communication", 6429, "Female", 2002, "Intermediate_2"),
    @("French", 809, "Female", 2002, "Intermediate_2"),
    @("Mathematics", 6113, "Female", 2002, "Intermediate_2"),
    @("Biology", 2351, "Female", 2002, "Intermediate_2"),
    @("Chemistry", 678, "Female", 2002, "Intermediate_2"),
    @("Physics", 380, "Female", 2002, "Intermediate_2"),
    @("Computing", 603, "Female", 2002, "Intermediate_2")
)

$row = 9
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $row++
}
